$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-7 down to 5-8
$ws.Rows.Item(4).Insert()

# Update S2 and S3 odds values (Odd_Over05_HT column)
$ws.Cells.Item(2, 19).Value = 1.47
$ws.Cells.Item(3, 19).Value = 1.58

# Ensure text columns (A-F) stay as text, not auto-converted (e.g. dates/times)
$ws.Range("A4:F8").NumberFormat = "@"

# Row 4: COLOMBIA - PRIMERA A - Atl. Nacional vs Ind. Medellin
$ws.Cells.Item(4, 1).Value = 'Ic6tHG75'
$ws.Cells.Item(4, 2).Value = '11/11/2024'
$ws.Cells.Item(4, 3).Value = '18:10'
$ws.Cells.Item(4, 4).Value = 'COLOMBIA - PRIMERA A'
$ws.Cells.Item(4, 5).Value = 'Atl. Nacional'
$ws.Cells.Item(4, 6).Value = 'Ind. Medellin'
$ws.Cells.Item(4, 7).Value = 2.25
$ws.Cells.Item(4, 8).Value = 3.4
$ws.Cells.Item(4, 9).Value = 3.1
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(4, 11).Value = 2.05
$ws.Cells.Item(4, 12).Value = 3.75
$ws.Cells.Item(4, 13).Value = 1.07
$ws.Cells.Item(4, 14).Value = 9
$ws.Cells.Item(4, 15).Value = 1.36
$ws.Cells.Item(4, 16).Value = 3
$ws.Cells.Item(4, 17).Value = 2.15
$ws.Cells.Item(4, 18).Value = 1.67
$ws.Cells.Item(4, 19).Value = 1.44
$ws.Cells.Item(4, 20).Value = 2.63
$ws.Cells.Item(4, 21).Value = 1.83
$ws.Cells.Item(4, 22).Value = 1.83
$ws.Cells.Item(4, 23).Value = 7
$ws.Cells.Item(4, 24).Value = 10
$ws.Cells.Item(4, 25).Value = 9.5
$ws.Cells.Item(4, 26).Value = 21
$ws.Cells.Item(4, 27).Value = 19
$ws.Cells.Item(4, 28).Value = 34
$ws.Cells.Item(4, 29).Value = 8.5
$ws.Cells.Item(4, 30).Value = 6.5
$ws.Cells.Item(4, 31).Value = 17
$ws.Cells.Item(4, 32).Value = 51
$ws.Cells.Item(4, 33).Value = 351
$ws.Cells.Item(4, 34).Value = 8.5
$ws.Cells.Item(4, 35).Value = 15
$ws.Cells.Item(4, 36).Value = 12
$ws.Cells.Item(4, 37).Value = 34
$ws.Cells.Item(4, 38).Value = 26
$ws.Cells.Item(4, 39).Value = 41
$ws.Cells.Item(4, 40).Value = 4.33
$ws.Cells.Item(4, 41).Value = 13
$ws.Cells.Item(4, 42).Value = 23
$ws.Cells.Item(4, 43).Value = 41
$ws.Cells.Item(4, 44).Value = 67
$ws.Cells.Item(4, 45).Value = 201
$ws.Cells.Item(4, 46).Value = 2.63
$ws.Cells.Item(4, 47).Value = 8.5
$ws.Cells.Item(4, 48).Value = 51
$ws.Cells.Item(4, 49).Value = 5
$ws.Cells.Item(4, 50).Value = 19
$ws.Cells.Item(4, 51).Value = 29
$ws.Cells.Item(4, 52).Value = 67
$ws.Cells.Item(4, 53).Value = 81
$ws.Cells.Item(4, 54).Value = 251
$ws.Cells.Item(4, 55).Value = 126
$ws.Cells.Item(4, 56).Value = 126

# Row 5: DENMARK - 1ST DIVISION - Hobro vs Hillerod
$ws.Cells.Item(5, 1).Value = 'G6XVkCH8'
$ws.Cells.Item(5, 2).Value = '11/11/2024'
$ws.Cells.Item(5, 3).Value = '15:00'
$ws.Cells.Item(5, 4).Value = 'DENMARK - 1ST DIVISION'
$ws.Cells.Item(5, 5).Value = 'Hobro'
$ws.Cells.Item(5, 6).Value = 'Hillerod'
$ws.Cells.Item(5, 7).Value = 2.5
$ws.Cells.Item(5, 8).Value = 3.6
$ws.Cells.Item(5, 9).Value = 2.5
$ws.Cells.Item(5, 10).Value = 3.1
$ws.Cells.Item(5, 11).Value = 2.25
$ws.Cells.Item(5, 12).Value = 3.1
$ws.Cells.Item(5, 13).Value = 1.03
$ws.Cells.Item(5, 14).Value = 15
$ws.Cells.Item(5, 15).Value = 1.2
$ws.Cells.Item(5, 16).Value = 4.33
$ws.Cells.Item(5, 17).Value = 1.67
$ws.Cells.Item(5, 18).Value = 2.15
$ws.Cells.Item(5, 19).Value = 1.33
$ws.Cells.Item(5, 20).Value = 3.25
$ws.Cells.Item(5, 21).Value = 1.57
$ws.Cells.Item(5, 22).Value = 2.25
$ws.Cells.Item(5, 23).Value = 11
$ws.Cells.Item(5, 24).Value = 15
$ws.Cells.Item(5, 25).Value = 10
$ws.Cells.Item(5, 26).Value = 26
$ws.Cells.Item(5, 27).Value = 19
$ws.Cells.Item(5, 28).Value = 23
$ws.Cells.Item(5, 29).Value = 13
$ws.Cells.Item(5, 30).Value = 7
$ws.Cells.Item(5, 31).Value = 12
$ws.Cells.Item(5, 32).Value = 41
$ws.Cells.Item(5, 33).Value = 151
$ws.Cells.Item(5, 34).Value = 11
$ws.Cells.Item(5, 35).Value = 15
$ws.Cells.Item(5, 36).Value = 10
$ws.Cells.Item(5, 37).Value = 26
$ws.Cells.Item(5, 38).Value = 19
$ws.Cells.Item(5, 39).Value = 23
$ws.Cells.Item(5, 40).Value = 4.75
$ws.Cells.Item(5, 41).Value = 13
$ws.Cells.Item(5, 42).Value = 21
$ws.Cells.Item(5, 43).Value = 41
$ws.Cells.Item(5, 44).Value = 51
$ws.Cells.Item(5, 45).Value = 126
$ws.Cells.Item(5, 46).Value = 3.25
$ws.Cells.Item(5, 47).Value = 7.5
$ws.Cells.Item(5, 48).Value = 41
$ws.Cells.Item(5, 49).Value = 4.75
$ws.Cells.Item(5, 50).Value = 13
$ws.Cells.Item(5, 51).Value = 21
$ws.Cells.Item(5, 52).Value = 41
$ws.Cells.Item(5, 53).Value = 51
$ws.Cells.Item(5, 54).Value = 126
$ws.Cells.Item(5, 55).Value = 351
$ws.Cells.Item(5, 56).Value = 126

# Row 6: ECUADOR - LIGA PRO - Imbabura vs Dep. Cuenca
$ws.Cells.Item(6, 1).Value = 'K8FqZP1I'
$ws.Cells.Item(6, 2).Value = '11/11/2024'
$ws.Cells.Item(6, 3).Value = '17:00'
$ws.Cells.Item(6, 4).Value = 'ECUADOR - LIGA PRO'
$ws.Cells.Item(6, 5).Value = 'Imbabura'
$ws.Cells.Item(6, 6).Value = 'Dep. Cuenca'
$ws.Cells.Item(6, 7).Value = 2.8
$ws.Cells.Item(6, 8).Value = 3.5
$ws.Cells.Item(6, 9).Value = 2.3
$ws.Cells.Item(6, 10).Value = 3.5
$ws.Cells.Item(6, 11).Value = 2.2
$ws.Cells.Item(6, 12).Value = 3
$ws.Cells.Item(6, 13).Value = 1.05
$ws.Cells.Item(6, 14).Value = 11
$ws.Cells.Item(6, 15).Value = 1.25
$ws.Cells.Item(6, 16).Value = 3.75
$ws.Cells.Item(6, 17).Value = 1.85
$ws.Cells.Item(6, 18).Value = 1.95
$ws.Cells.Item(6, 19).Value = 1.4
$ws.Cells.Item(6, 20).Value = 2.75
$ws.Cells.Item(6, 21).Value = 1.73
$ws.Cells.Item(6, 22).Value = 2
$ws.Cells.Item(6, 23).Value = 9.5
$ws.Cells.Item(6, 24).Value = 15
$ws.Cells.Item(6, 25).Value = 11
$ws.Cells.Item(6, 26).Value = 29
$ws.Cells.Item(6, 27).Value = 23
$ws.Cells.Item(6, 28).Value = 29
$ws.Cells.Item(6, 29).Value = 11
$ws.Cells.Item(6, 30).Value = 6.5
$ws.Cells.Item(6, 31).Value = 15
$ws.Cells.Item(6, 32).Value = 51
$ws.Cells.Item(6, 33).Value = 201
$ws.Cells.Item(6, 34).Value = 8.5
$ws.Cells.Item(6, 35).Value = 12
$ws.Cells.Item(6, 36).Value = 9.5
$ws.Cells.Item(6, 37).Value = 21
$ws.Cells.Item(6, 38).Value = 19
$ws.Cells.Item(6, 39).Value = 29
$ws.Cells.Item(6, 40).Value = 5
$ws.Cells.Item(6, 41).Value = 17
$ws.Cells.Item(6, 42).Value = 23
$ws.Cells.Item(6, 43).Value = 51
$ws.Cells.Item(6, 44).Value = 67
$ws.Cells.Item(6, 45).Value = 151
$ws.Cells.Item(6, 46).Value = 2.75
$ws.Cells.Item(6, 47).Value = 7.5
$ws.Cells.Item(6, 48).Value = 51
$ws.Cells.Item(6, 49).Value = 4.5
$ws.Cells.Item(6, 50).Value = 13
$ws.Cells.Item(6, 51).Value = 21
$ws.Cells.Item(6, 52).Value = 41
$ws.Cells.Item(6, 53).Value = 51
$ws.Cells.Item(6, 54).Value = 151
$ws.Cells.Item(6, 55).Value = 51
$ws.Cells.Item(6, 56).Value = 51

# Row 7: ROMANIA - LIGA 1 - Otelul vs CFR Cluj
$ws.Cells.Item(7, 1).Value = 'Obk4TfWF'
$ws.Cells.Item(7, 2).Value = '11/11/2024'
$ws.Cells.Item(7, 3).Value = '16:30'
$ws.Cells.Item(7, 4).Value = 'ROMANIA - LIGA 1'
$ws.Cells.Item(7, 5).Value = 'Otelul'
$ws.Cells.Item(7, 6).Value = 'CFR Cluj'
$ws.Cells.Item(7, 7).Value = 4
$ws.Cells.Item(7, 8).Value = 2.88
$ws.Cells.Item(7, 9).Value = 2.15
$ws.Cells.Item(7, 10).Value = 4.5
$ws.Cells.Item(7, 11).Value = 1.91
$ws.Cells.Item(7, 12).Value = 3
$ws.Cells.Item(7, 13).Value = 1.11
$ws.Cells.Item(7, 14).Value = 6.5
$ws.Cells.Item(7, 15).Value = 1.5
$ws.Cells.Item(7, 16).Value = 2.5
$ws.Cells.Item(7, 17).Value = 2.6
$ws.Cells.Item(7, 18).Value = 1.48
$ws.Cells.Item(7, 19).Value = 1.57
$ws.Cells.Item(7, 20).Value = 2.25
$ws.Cells.Item(7, 21).Value = 2.2
$ws.Cells.Item(7, 22).Value = 1.62
$ws.Cells.Item(7, 23).Value = 8.5
$ws.Cells.Item(7, 24).Value = 17
$ws.Cells.Item(7, 25).Value = 15
$ws.Cells.Item(7, 26).Value = 41
$ws.Cells.Item(7, 27).Value = 41
$ws.Cells.Item(7, 28).Value = 51
$ws.Cells.Item(7, 29).Value = 6
$ws.Cells.Item(7, 30).Value = 6
$ws.Cells.Item(7, 31).Value = 19
$ws.Cells.Item(7, 32).Value = 67
$ws.Cells.Item(7, 33).Value = 351
$ws.Cells.Item(7, 34).Value = 6
$ws.Cells.Item(7, 35).Value = 9
$ws.Cells.Item(7, 36).Value = 10
$ws.Cells.Item(7, 37).Value = 19
$ws.Cells.Item(7, 38).Value = 21
$ws.Cells.Item(7, 39).Value = 41
$ws.Cells.Item(7, 40).Value = 5.5
$ws.Cells.Item(7, 41).Value = 23
$ws.Cells.Item(7, 42).Value = 41
$ws.Cells.Item(7, 43).Value = 81
$ws.Cells.Item(7, 44).Value = 126
$ws.Cells.Item(7, 45).Value = 351
$ws.Cells.Item(7, 46).Value = 2.25
$ws.Cells.Item(7, 47).Value = 9.5
$ws.Cells.Item(7, 48).Value = 81
$ws.Cells.Item(7, 49).Value = 4
$ws.Cells.Item(7, 50).Value = 13
$ws.Cells.Item(7, 51).Value = 29
$ws.Cells.Item(7, 52).Value = 41
$ws.Cells.Item(7, 53).Value = 81
$ws.Cells.Item(7, 54).Value = 301
$ws.Cells.Item(7, 55).Value = 51
$ws.Cells.Item(7, 56).Value = 51

# Row 8: SPAIN - LALIGA2 - Dep. La Coruna vs Eibar
$ws.Cells.Item(8, 1).Value = 'Wj3sKyQp'
$ws.Cells.Item(8, 2).Value = '11/11/2024'
$ws.Cells.Item(8, 3).Value = '16:30'
$ws.Cells.Item(8, 4).Value = 'SPAIN - LALIGA2'
$ws.Cells.Item(8, 5).Value = 'Dep. La Coruna'
$ws.Cells.Item(8, 6).Value = 'Eibar'
$ws.Cells.Item(8, 7).Value = 1.73
$ws.Cells.Item(8, 8).Value = 3.5
$ws.Cells.Item(8, 9).Value = 5
$ws.Cells.Item(8, 10).Value = 2.38
$ws.Cells.Item(8, 11).Value = 2.2
$ws.Cells.Item(8, 12).Value = 5
$ws.Cells.Item(8, 13).Value = 1.06
$ws.Cells.Item(8, 14).Value = 9.5
$ws.Cells.Item(8, 15).Value = 1.3
$ws.Cells.Item(8, 16).Value = 3.4
$ws.Cells.Item(8, 17).Value = 2.03
$ws.Cells.Item(8, 18).Value = 1.83
$ws.Cells.Item(8, 19).Value = 1.4
$ws.Cells.Item(8, 20).Value = 2.75
$ws.Cells.Item(8, 21).Value = 1.83
$ws.Cells.Item(8, 22).Value = 1.83
$ws.Cells.Item(8, 23).Value = 6.5
$ws.Cells.Item(8, 24).Value = 8
$ws.Cells.Item(8, 25).Value = 8.5
$ws.Cells.Item(8, 26).Value = 13
$ws.Cells.Item(8, 27).Value = 15
$ws.Cells.Item(8, 28).Value = 29
$ws.Cells.Item(8, 29).Value = 9.5
$ws.Cells.Item(8, 30).Value = 7
$ws.Cells.Item(8, 31).Value = 17
$ws.Cells.Item(8, 32).Value = 51
$ws.Cells.Item(8, 33).Value = 301
$ws.Cells.Item(8, 34).Value = 13
$ws.Cells.Item(8, 35).Value = 23
$ws.Cells.Item(8, 36).Value = 17
$ws.Cells.Item(8, 37).Value = 51
$ws.Cells.Item(8, 38).Value = 41
$ws.Cells.Item(8, 39).Value = 41
$ws.Cells.Item(8, 40).Value = 3.6
$ws.Cells.Item(8, 41).Value = 9
$ws.Cells.Item(8, 42).Value = 21
$ws.Cells.Item(8, 43).Value = 29
$ws.Cells.Item(8, 44).Value = 51
$ws.Cells.Item(8, 45).Value = 151
$ws.Cells.Item(8, 46).Value = 2.75
$ws.Cells.Item(8, 47).Value = 8.5
$ws.Cells.Item(8, 48).Value = 51
$ws.Cells.Item(8, 49).Value = 6.5
$ws.Cells.Item(8, 50).Value = 26
$ws.Cells.Item(8, 51).Value = 34
$ws.Cells.Item(8, 52).Value = 101
$ws.Cells.Item(8, 53).Value = 126
$ws.Cells.Item(8, 54).Value = 251
$ws.Cells.Item(8, 55).Value = 81
$ws.Cells.Item(8, 56).Value = 81
